$d = $word.ActiveDocument

$replacements = @(
    @("91÷7=", "52÷4="),
    @("52÷9=", "35÷2="),
    @("46÷5=", "95÷5="),
    @("19÷6=", "94÷5="),
    @("64÷9=", "89÷8="),
    @("98÷6=", "92÷4="),
    @("52÷7=", "91÷6="),
    @("73÷8=", "31÷2="),
    @("90÷6=", "57÷6="),
    @("60÷4=", "10÷9="),
    @("49÷7=", "54÷5="),
    @("18÷9=", "85÷7="),
    @("27÷5=", "54÷7="),
    @("47÷9=", "56÷9="),
    @("99÷2=", "55÷2="),
    @("83÷6=", "62÷7="),
    @("66÷4=", "74÷8="),
    @("73÷4=", "51÷3="),
    @("84÷2=", "28÷2="),
    @("95÷4=", "80÷2="),
    @("33÷7=", "69÷8="),
    @("13÷9=", "71÷8="),
    @("61÷5=", "54÷9="),
    @("22÷6=", "88÷8="),
    @("99÷6=", "16÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
